$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update C2 to the new IP address
$ws.Range("C2").Value = "127.0.0.1"

# E2 keeps its IP value, but now adopts the same (text) style as C2
$ws.Range("E2").Value = "192.168.0.24"
$ws.Range("E2").NumberFormat = $ws.Range("C2").NumberFormat

# Move the active selection from C2 to E2
$ws.Range("E2").Select()
